# Insert 2 new data rows into the Ciruela price table right before the
# current row 303 (which holds the 2022-02-24 "Angeleno" records). This
# pushes every existing row down by two (old 303..351 -> new 305..353)
# and creates space for a brand-new "Friar" variety record dated 44984
# (2023-02-08) at new rows 303..304, matching the published diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 303; Excel shifts the
# existing rows (and the sheet dimension) down automatically.
$ws.Rows(303).Resize(2).Insert()

# Row 303: Friar / Primera
$ws.Cells.Item(303, 1).Value2  = 8
$ws.Cells.Item(303, 2).Value   = "Terminal La Palmera de La Serena"
$ws.Cells.Item(303, 3).Value   = "Coquimbo"
$ws.Cells.Item(303, 4).Value2  = 44984
$ws.Cells.Item(303, 5).Value2  = 4
$ws.Cells.Item(303, 6).Value   = "Fruta"
$ws.Cells.Item(303, 7).Value2  = 100103
$ws.Cells.Item(303, 8).Value   = "Frutos de hueso (carozo)"
$ws.Cells.Item(303, 9).Value2  = 100103002
$ws.Cells.Item(303, 10).Value  = "Ciruela"
$ws.Cells.Item(303, 11).Value  = "Friar"
$ws.Cells.Item(303, 12).Value  = "Primera"
$ws.Cells.Item(303, 13).Value2 = 20
$ws.Cells.Item(303, 14).Value2 = 210000
$ws.Cells.Item(303, 15).Value2 = 220000
$ws.Cells.Item(303, 16).Value2 = 215000
$ws.Cells.Item(303, 17).Value  = "`$/bins (450 kilos)"
$ws.Cells.Item(303, 18).Value  = "Región Metropolitana"
$ws.Cells.Item(303, 19).Value2 = 478
$ws.Cells.Item(303, 20).Value2 = 450

# Row 304: Friar / Segunda
$ws.Cells.Item(304, 1).Value2  = 8
$ws.Cells.Item(304, 2).Value   = "Terminal La Palmera de La Serena"
$ws.Cells.Item(304, 3).Value   = "Coquimbo"
$ws.Cells.Item(304, 4).Value2  = 44984
$ws.Cells.Item(304, 5).Value2  = 4
$ws.Cells.Item(304, 6).Value   = "Fruta"
$ws.Cells.Item(304, 7).Value2  = 100103
$ws.Cells.Item(304, 8).Value   = "Frutos de hueso (carozo)"
$ws.Cells.Item(304, 9).Value2  = 100103002
$ws.Cells.Item(304, 10).Value  = "Ciruela"
$ws.Cells.Item(304, 11).Value  = "Friar"
$ws.Cells.Item(304, 12).Value  = "Segunda"
$ws.Cells.Item(304, 13).Value2 = 20
$ws.Cells.Item(304, 14).Value2 = 170000
$ws.Cells.Item(304, 15).Value2 = 180000
$ws.Cells.Item(304, 16).Value2 = 175000
$ws.Cells.Item(304, 17).Value  = "`$/bins (450 kilos)"
$ws.Cells.Item(304, 18).Value  = "Región Metropolitana"
$ws.Cells.Item(304, 19).Value2 = 389
$ws.Cells.Item(304, 20).Value2 = 450

# Keep the date formatting consistent with the rest of column D.
$ws.Range("D303:D304").NumberFormat = $ws.Range("D305").NumberFormat
